# Public_lpAcctData.xlsx / Sheet1
#
# The sheet gains a header row and two new trailing columns (UserStatus, UserName);
# the GroupCode date column and the UserEmail column keep their original data (just
# shifted down one row), and the column widths get rearranged: the old "dates" column
# loses its custom width, the old "emails" column width (32.7109375) ends up on the new
# 4th column, and the email column itself gets a brand-new custom width.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Make room for a header row - existing rows 1-4 become rows 2-5.
$ws.Rows.Item(1).Insert()

# 2) Make room for a new column at B so the existing columns (and their exact widths)
#    shift out of the way: old B (dates, width 14.7109375) -> C, old C (emails, width
#    32.7109375) -> D.
$ws.Columns.Item(2).Insert()

# 3) Move the shifted data back into place using a values-only paste, so the
#    destination columns don't inherit the source columns' widths.
$ws.Range("C2:C5").Copy()
$ws.Range("B2:B5").PasteSpecial(-4163)

$ws.Range("D2:D5").Copy()
$ws.Range("C2:C5").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# 4) New header row.
$ws.Cells.Item(1, 2).Value = "GroupCode"
$ws.Cells.Item(1, 3).Value = "UserEmail"
$ws.Cells.Item(1, 4).Value = "UserStatus"
$ws.Cells.Item(1, 5).Value = "UserName"

# 5) Refresh the GroupCode dates.
for ($r = 2; $r -le 5; $r++) {
    $ws.Cells.Item($r, 2).Value = 38174
}

# 6) New UserStatus / UserName data.
$ws.Cells.Item(2, 4).Value = "Active"
$ws.Cells.Item(2, 5).Value = "lpuser-5-23-2012-55611"

$ws.Cells.Item(3, 4).Value = "Active"
$ws.Cells.Item(3, 5).Value = "lpuser-5-23-2012-55645"

$ws.Cells.Item(4, 4).Value = "Inactive"
$ws.Cells.Item(4, 5).Value = "lpuser-5-23-2012-55681"

$ws.Cells.Item(5, 4).Value = "Active"
$ws.Cells.Item(5, 5).Value = "lpuser-5-23-2012-55713"

# 7) Widen the UserEmail column to its new size (column D already inherited the exact
#    32.7109375 width from the old emails column during the shift in step 2).
$ws.Columns.Item(3).ColumnWidth = 28

# 8) Match the saved selection.
$ws.Range("E11:E12").Select()
